$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 31 becomes a DIALOGUE line carrying the Monologue text that used to live on row 32
$ws.Range("B31").Value = "DIALOGUE"
$ws.Range("C31").Value = "Monologue"
$ws.Range("D31").Value = "예진은 낯빛이 어두워지며 뒤돌아 자신의 자리로 향한다."

# Row 32 now just holds the HIDE_CHAR command that used to be on row 31, columns C/D fully cleared
$ws.Range("C32:D32").Clear()
$ws.Range("B32").Value = "HIDE_CHAR"

# Update the active selection to match the authored state
$ws.Range("D7").Select()
